# Update table 2 and 3
# Refresh the "Ratio MoMo/COVID-19 Deaths" column (E) with more precise
# (2-decimal) ratio values instead of the previously rounded 1-decimal
# values, and switch the number format used for that column from the
# custom "0.0" format to the built-in "0.00" format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 ("all" subtotal for 2021) already used a dedicated number
# format (0.0, the only cell on the sheet using it). Updating it first
# lets the engine rewrite that format slot in place as "0.00" instead
# of allocating a brand-new style slot.
$ws.Range("E12").Value = 0.54
$ws.Range("E12").NumberFormat = "0.00"

# Row 7 is the "all" subtotal for 2020; it now also picks up the same
# 0.00 number format (previously it used the plain General format).
$ws.Range("E7").Value = 0.5
$ws.Range("E7").NumberFormat = "0.00"

# --- 2020 block (rows 3-6) ---
$ws.Range("E3").Value = 0.15
$ws.Range("E4").Value = 0.58
$ws.Range("E5").Value = 0.59
$ws.Range("E6").Value = 0.43

# --- 2021 block (rows 8-11) ---
$ws.Range("E8").Value = 0.27
$ws.Range("E9").Value = 0.65
$ws.Range("E10").Value = 0.62
$ws.Range("E11").Value = 0.45
